$d = $word.ActiveDocument

# The row describing "whenCaptured" has its Type cell containing just "date".
# We need to append a new, separate run with text "time" (same formatting:
# sz 28 / 14pt, szCs 28, lang en-US) right after the existing "date" run,
# inside the same paragraph/cell, so the cell reads "datetime" made up of
# two runs: "date" (untouched) + "time" (new).
#
# A plain Range.InsertAfter() right next to text that already has identical
# character formatting gets silently coalesced into the existing run when
# the document is serialized (Word/this host normalizes adjacent runs that
# share formatting). To force the new text to remain its own <w:r>, we
# briefly turn on revision tracking while inserting the text (which always
# produces a distinct run for the inserted text) and then accept that one
# revision - leaving two plain (non-tracked) runs behind.

$table = $d.Tables.Item(1)

# Row 7 / column 2 (1-based) is the "whenCaptured" row's Type cell, whose
# text is just "date" - but look it up by content instead of a hard-coded
# index so the script keeps working even if the table shape changes.
$cell = $null
for ($r = 1; $r -le $table.Rows.Count; $r++) {
    $candidate = $table.Cell($r, 2)
    $t = $candidate.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "date") {
        $cell = $candidate
        break
    }
}
if ($cell -eq $null) {
    $cell = $table.Cell(7, 2)
}

$cellRange = $cell.Range
# Cell.Range includes the trailing end-of-cell marker; stop just before it.
$insertRange = $cellRange.Duplicate
$insertRange.Collapse(0)      # wdCollapseEnd
$insertRange.MoveEnd(1, -1)   # back up over the end-of-cell mark

$wasTracking = $d.TrackRevisions
$d.TrackRevisions = $true
$insertRange.InsertAfter("time")
$d.TrackRevisions = $wasTracking

# Accept the single revision we just created so the inserted text becomes
# normal (non-tracked) content, kept as its own run.
if ($d.Revisions.Count -gt 0) {
    $d.Revisions.Item($d.Revisions.Count).Accept()
}
